$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 2 (item id 5489)
$ws.Range("H2").Value = 228.71428
$ws.Range("I2").Value = 125.25
$ws.Range("J2").Value = 366.66666
$ws.Range("K2").Value = 125.25
$ws.Range("L2").Value = 366.66666
$ws.Range("M2").Value = -12.25
$ws.Range("N2").Value = -592.66666
# Row 12 (item id 5515)
$ws.Range("H12").Value = 324.66666
$ws.Range("J12").Value = 339.6
$ws.Range("L12").Value = 339.6
$ws.Range("N12").Value = -679.6
# Row 18 (item id 5471)
$ws.Range("H18").Value = 9299.429
$ws.Range("I18").Value = 1239.6
$ws.Range("J18").Value = 13777.111
$ws.Range("K18").Value = 1239.6
$ws.Range("L18").Value = 13777.111
$ws.Range("M18").Value = -955.5999999999999
$ws.Range("N18").Value = -14345.111
# Row 87 (item id 10651)
$ws.Range("H87").Value = 70181.75
$ws.Range("J87").Value = 91490.8
$ws.Range("L87").Value = 91490.8
$ws.Range("N87").Value = -93986.8
# Row 90 (item id 10651)
$ws.Range("H90").Value = 70181.75
$ws.Range("J90").Value = 91490.8
$ws.Range("L90").Value = 274472.4
$ws.Range("N90").Value = -286952.4
# Row 106 (item id 19903)
$ws.Range("H106").Value = 2632.5
$ws.Range("I106").Value = 2159
$ws.Range("K106").Value = 2159
$ws.Range("M106").Value = -1528
# Row 112 (item id 27960)
$ws.Range("H112").Value = 1620.0377
$ws.Range("I112").Value = 844.6667
$ws.Range("J112").Value = 1666.56
$ws.Range("K112").Value = 2534.0001
$ws.Range("L112").Value = 4999.68
$ws.Range("M112").Value = -1426.0001
$ws.Range("N112").Value = -7215.68
# Row 113 (item id 27775)
$ws.Range("H113").Value = 3508.182
$ws.Range("I113").Value = 3310
$ws.Range("K113").Value = 3310
$ws.Range("M113").Value = -56
# Row 116 (item id 27778)
$ws.Range("H116").Value = 57147.93
$ws.Range("I116").Value = 97359.28999999999
$ws.Range("J116").Value = 16936.572
$ws.Range("K116").Value = 97359.28999999999
$ws.Range("L116").Value = 16936.572
$ws.Range("M116").Value = -93917.28999999999
$ws.Range("N116").Value = -23820.572
# Row 138 (item id 44169)
$ws.Range("H138").Value = 3073.18
$ws.Range("I138").Value = 2719.2144
$ws.Range("K138").Value = 8157.6432
$ws.Range("M138").Value = -3017.6432

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2 (item id 27713)
$ws.Range("H2").Value = 2227.6333
$ws.Range("I2").Value = 2049.1738
$ws.Range("K2").Value = 2049.1738
$ws.Range("M2").Value = -1936.1738
# Row 32 (item id 44147)
$ws.Range("H32").Value = 1668639.6
$ws.Range("I32").Value = 1685491.5
$ws.Range("K32").Value = 1685491.5
$ws.Range("M32").Value = -1685204.5
# Row 45 (item id 27714)
$ws.Range("H45").Value = 5288.5386
$ws.Range("I45").Value = 5816.8887
$ws.Range("K45").Value = 5816.8887
$ws.Range("M45").Value = -5439.8887
# Row 51 (item id 3858)
$ws.Range("H51").Value = 40000
$ws.Range("J51").Value = 40000
$ws.Range("L51").Value = 40000
$ws.Range("N51").Value = -41512
# Row 61 (item id 43999)
$ws.Range("H61").Value = 857734.1
$ws.Range("I61").Value = 1012249.44
$ws.Range("K61").Value = 1012249.44
$ws.Range("M61").Value = -1012037.44
# Row 74 (item id 44000)
$ws.Range("H74").Value = 3909131.5
$ws.Range("I74").Value = 4631643
$ws.Range("K74").Value = 4631643
$ws.Range("M74").Value = -4630769
# Row 77 (item id 44000)
$ws.Range("H77").Value = 3909131.5
$ws.Range("I77").Value = 4631643
$ws.Range("K77").Value = 23158215
$ws.Range("M77").Value = -23153847
# Row 116 (item id 27713)
$ws.Range("H116").Value = 2227.6333
$ws.Range("I116").Value = 2049.1738
$ws.Range("K116").Value = 2049.1738
$ws.Range("M116").Value = 244.8262
# Row 122 (item id 36168)
$ws.Range("H122").Value = 3279.9644
$ws.Range("I122").Value = 2906.0435
$ws.Range("K122").Value = 8718.130500000001
$ws.Range("M122").Value = -6268.130500000001
# Row 128 (item id 34570)
$ws.Range("H128").Value = 90000
$ws.Range("J128").Value = 90000
$ws.Range("L128").Value = 90000
$ws.Range("N128").Value = -99960
# Row 136 (item id 43999)
$ws.Range("H136").Value = 857734.1
$ws.Range("I136").Value = 1012249.44
$ws.Range("K136").Value = 3036748.32
$ws.Range("M136").Value = -3034198.32

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3 (item id 27713)
$ws.Range("H3").Value = 2227.6333
$ws.Range("I3").Value = 2049.1738
$ws.Range("K3").Value = 2049.1738
$ws.Range("M3").Value = -1935.1738
# Row 20 (item id 14149)
$ws.Range("H20").Value = 1864.1364
$ws.Range("I20").Value = 2028.5
$ws.Range("J20").Value = 1666.9
$ws.Range("K20").Value = 2028.5
$ws.Range("L20").Value = 1666.9
$ws.Range("M20").Value = -1781.5
$ws.Range("N20").Value = -2160.9
# Row 22 (item id 5092)
$ws.Range("H22").Value = 824.75
$ws.Range("I22").Value = 824.75
$ws.Range("K22").Value = 824.75
$ws.Range("M22").Value = -651.75

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 2 (item id 1820)
$ws.Range("H2").Value = 144
$ws.Range("I2").Value = 144
$ws.Range("K2").Value = 144
$ws.Range("M2").Value = -31
# Row 14 (item id 1998)
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
# Row 22 (item id 5367)
$ws.Range("H22").Value = 948.6923
$ws.Range("I22").Value = 858.63635
$ws.Range("K22").Value = 858.63635
$ws.Range("M22").Value = -508.63635
# Row 31 (item id 44023)
$ws.Range("H31").Value = 3858.524
$ws.Range("I31").Value = 2154.0967
$ws.Range("J31").Value = 5509.6875
$ws.Range("K31").Value = 2154.0967
$ws.Range("L31").Value = 5509.6875
$ws.Range("M31").Value = -1859.0967
$ws.Range("N31").Value = -6099.6875
# Row 34 (item id 44023)
$ws.Range("H34").Value = 3858.524
$ws.Range("I34").Value = 2154.0967
$ws.Range("J34").Value = 5509.6875
$ws.Range("K34").Value = 2154.0967
$ws.Range("L34").Value = 5509.6875
$ws.Range("M34").Value = -1952.0967
$ws.Range("N34").Value = -5913.6875
# Row 48 (item id 3870)
$ws.Range("H48").Value = 58999
$ws.Range("J48").Value = 58999
$ws.Range("L48").Value = 58999
$ws.Range("N48").Value = -59951
# Row 68 (item id 10611)
$ws.Range("H68").Value = 70996.86
$ws.Range("I68").Value = 49000
$ws.Range("J68").Value = 72688.92
$ws.Range("K68").Value = 49000
$ws.Range("L68").Value = 72688.92
$ws.Range("M68").Value = -48251
$ws.Range("N68").Value = -74186.92
# Row 71 (item id 10611)
$ws.Range("H71").Value = 70996.86
$ws.Range("I71").Value = 49000
$ws.Range("J71").Value = 72688.92
$ws.Range("K71").Value = 147000
$ws.Range("L71").Value = 218066.76
$ws.Range("M71").Value = -143256
$ws.Range("N71").Value = -225554.76
# Row 132 (item id 44019)
$ws.Range("H132").Value = 5327950
$ws.Range("I132").Value = 9372.762000000001
$ws.Range("K132").Value = 28118.286
$ws.Range("M132").Value = -25588.286

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 11 (item id 4745)
$ws.Range("H11").Value = 175.5
$ws.Range("I11").Value = 175.5
$ws.Range("K11").Value = 526.5
$ws.Range("M11").Value = -386.5
# Row 37 (item id 9516)
$ws.Range("H37").Value = 89991
$ws.Range("J37").Value = 89991
$ws.Range("L37").Value = 269973
$ws.Range("N37").Value = -270197
# Row 75 (item id 12863)
$ws.Range("H75").Value = 6371.8184
$ws.Range("I75").Value = 1049
$ws.Range("J75").Value = 6904.1
$ws.Range("K75").Value = 3147
$ws.Range("L75").Value = 20712.3
$ws.Range("M75").Value = -2149
$ws.Range("N75").Value = -22708.3
# Row 78 (item id 12863)
$ws.Range("H78").Value = 6371.8184
$ws.Range("I78").Value = 1049
$ws.Range("J78").Value = 6904.1
$ws.Range("K78").Value = 9441
$ws.Range("L78").Value = 62136.9
$ws.Range("M78").Value = -4449
$ws.Range("N78").Value = -72120.89999999999

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 49 (item id 4232)
$ws.Range("H49").Value = 45000
$ws.Range("J49").Value = 45000
$ws.Range("L49").Value = 45000
$ws.Range("N49").Value = -45368
# Row 70 (item id 14146)
$ws.Range("H70").Value = 9932.333000000001
$ws.Range("I70").Value = 9932.333000000001
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 9932.333000000001
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -9662.333000000001
$ws.Range("N70").ClearContents()
# Row 73 (item id 14146)
$ws.Range("H73").Value = 9932.333000000001
$ws.Range("I73").Value = 9932.333000000001
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 9932.333000000001
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -8996.333000000001
$ws.Range("N73").ClearContents()

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 20 (item id 4308)
$ws.Range("H20").Value = 15006
$ws.Range("J20").Value = 15006
$ws.Range("L20").Value = 15006
$ws.Range("N20").Value = -15458
# Row 109 (item id 27209)
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
# Row 122 (item id 36247)
$ws.Range("H122").Value = 3485.2646
$ws.Range("I122").Value = 3255.1538
$ws.Range("K122").Value = 9765.4614
$ws.Range("M122").Value = -7315.4614

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 4 (item id 2996)
$ws.Range("H4").Value = 167607840
$ws.Range("J4").Value = 200029390
$ws.Range("L4").Value = 200029390
$ws.Range("N4").Value = -200029616
# Row 132 (item id 44029)
$ws.Range("H132").Value = 6101697
$ws.Range("I132").Value = 6942034.5
$ws.Range("K132").Value = 20826103.5
$ws.Range("M132").Value = -20823573.5
